$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 44007876
$ws.Range("J70").Value = 55563972
$ws.Range("L70").Value = 166691916
$ws.Range("N70").Value = -166692456
$ws.Range("H73").Value = 44007876
$ws.Range("J73").Value = 55563972
$ws.Range("L73").Value = 166691916
$ws.Range("N73").Value = -166693788
$ws.Range("H100").Value = 5508.2144
$ws.Range("I100").Value = 2334.818
$ws.Range("J100").Value = 7561.5884
$ws.Range("K100").Value = 2334.818
$ws.Range("L100").Value = 7561.5884
$ws.Range("M100").Value = -1793.818
$ws.Range("N100").Value = -8643.588400000001
$ws.Range("H111").Value = 1090.8334
$ws.Range("I111").Value = 1548.3334
$ws.Range("J111").Value = 633.3333
$ws.Range("K111").Value = 4645.0002
$ws.Range("L111").Value = 1899.9999
$ws.Range("M111").Value = -1578.0002
$ws.Range("N111").Value = -8033.9999
$ws.Range("H116").Value = 4987.3335
$ws.Range("I116").Value = 4650
$ws.Range("K116").Value = 4650
$ws.Range("M116").Value = -1208
$ws.Range("H132").Value = 1699.5927
$ws.Range("I132").Value = 1375.8096
$ws.Range("J132").Value = 2832.8333
$ws.Range("K132").Value = 4127.4288
$ws.Range("L132").Value = 8498.499899999999
$ws.Range("M132").Value = -1597.4288
$ws.Range("N132").Value = -13558.4999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2835.2
$ws.Range("I32").Value = 2775.0444
$ws.Range("K32").Value = 2775.0444
$ws.Range("M32").Value = -2488.0444
$ws.Range("H34").Value = 105074.62
$ws.Range("I34").Value = 68994.25
$ws.Range("K34").Value = 68994.25
$ws.Range("M34").Value = -68723.25
$ws.Range("H61").Value = 4659.409
$ws.Range("J61").Value = 1932.3334
$ws.Range("L61").Value = 1932.3334
$ws.Range("N61").Value = -2356.3334
$ws.Range("H82").Value = 112590.5
$ws.Range("J82").Value = 112590.5
$ws.Range("L82").Value = 112590.5
$ws.Range("N82").Value = -113312.5
$ws.Range("H85").Value = 112590.5
$ws.Range("J85").Value = 112590.5
$ws.Range("L85").Value = 112590.5
$ws.Range("N85").Value = -115086.5
$ws.Range("H105").Value = 60335
$ws.Range("I105").Value = 60335
$ws.Range("K105").Value = 60335
$ws.Range("M105").Value = -56841
$ws.Range("H136").Value = 4659.409
$ws.Range("J136").Value = 1932.3334
$ws.Range("L136").Value = 5797.0002
$ws.Range("N136").Value = -10897.0002

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 354.48386
$ws.Range("I22").Value = 354.48386
$ws.Range("K22").Value = 354.48386
$ws.Range("M22").Value = -181.48386
$ws.Range("H94").Value = 3618.524
$ws.Range("I94").Value = 3202.0588
$ws.Range("J94").Value = 5388.5
$ws.Range("K94").Value = 3202.0588
$ws.Range("L94").Value = 5388.5
$ws.Range("M94").Value = -2751.0588
$ws.Range("N94").Value = -6290.5
$ws.Range("H99").Value = 2986.8696
$ws.Range("I99").Value = 2748.5334
$ws.Range("J99").Value = 3433.75
$ws.Range("K99").Value = 2748.5334
$ws.Range("L99").Value = 3433.75
$ws.Range("M99").Value = -1250.5334
$ws.Range("N99").Value = -6429.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 320.125
$ws.Range("I7").Value = 70
$ws.Range("J7").Value = 403.5
$ws.Range("K7").Value = 70
$ws.Range("L7").Value = 403.5
$ws.Range("M7").Value = 43
$ws.Range("N7").Value = -629.5
$ws.Range("H42").Value = 2000
$ws.Range("I42").Value = 2000
$ws.Range("K42").Value = 2000
$ws.Range("M42").Value = -1407
$ws.Range("H45").Value = 15000
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H62").Value = 9131.154
$ws.Range("I62").Value = 9484.166999999999
$ws.Range("J62").Value = 8828.571
$ws.Range("K62").Value = 9484.166999999999
$ws.Range("L62").Value = 8828.571
$ws.Range("M62").Value = -8860.166999999999
$ws.Range("N62").Value = -10076.571
$ws.Range("H65").Value = 9131.154
$ws.Range("I65").Value = 9484.166999999999
$ws.Range("J65").Value = 8828.571
$ws.Range("K65").Value = 47420.835
$ws.Range("L65").Value = 44142.855
$ws.Range("M65").Value = -44300.835
$ws.Range("N65").Value = -50382.855
$ws.Range("H93").Value = 49933
$ws.Range("I93").Value = 49899.5
$ws.Range("J93").Value = 50000
$ws.Range("K93").Value = 49899.5
$ws.Range("L93").Value = 50000
$ws.Range("M93").Value = -48027.5
$ws.Range("N93").Value = -53744
$ws.Range("H103").Value = 11604.4
$ws.Range("I103").Value = 11604.4
$ws.Range("K103").Value = 11604.4
$ws.Range("M103").Value = -10432.4
$ws.Range("H134").Value = 3439.389
$ws.Range("I134").Value = 3649.3333
$ws.Range("J134").Value = 2389.6667
$ws.Range("K134").Value = 10947.9999
$ws.Range("L134").Value = 7169.000100000001
$ws.Range("M134").Value = -8412.999899999999
$ws.Range("N134").Value = -12239.0001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 761.55
$ws.Range("I121").Value = 621.125
$ws.Range("J121").Value = 855.1667
$ws.Range("K121").Value = 1863.375
$ws.Range("L121").Value = 2565.5001
$ws.Range("M121").Value = -553.375
$ws.Range("N121").Value = -5185.5001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9500
$ws.Range("H73").Value = 9500
$ws.Range("H80").Value = 6191.25
$ws.Range("I80").Value = 4098.3335
$ws.Range("K80").Value = 4098.3335
$ws.Range("M80").Value = -3100.3335
$ws.Range("H83").Value = 6191.25
$ws.Range("I83").Value = 4098.3335
$ws.Range("K83").Value = 20491.6675
$ws.Range("M83").Value = -15499.6675

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H40").Value = 6431.4062
$ws.Range("I40").Value = 4973.8667
$ws.Range("J40").Value = 7717.4707
$ws.Range("K40").Value = 4973.8667
$ws.Range("L40").Value = 7717.4707
$ws.Range("M40").Value = -4837.8667
$ws.Range("N40").Value = -7989.4707
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H104").Value = 54813.75
$ws.Range("J104").Value = 54813.75
$ws.Range("L104").Value = 54813.75
$ws.Range("N104").Value = -61801.75
$ws.Range("H132").Value = 3521.5
$ws.Range("I132").Value = 3465.8
$ws.Range("J132").Value = 3800
$ws.Range("K132").Value = 10397.4
$ws.Range("L132").Value = 11400
$ws.Range("M132").Value = -7867.400000000001
$ws.Range("N132").Value = -16460

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 17562.5
$ws.Range("I31").Value = 15714.286
$ws.Range("K31").Value = 15714.286
$ws.Range("M31").Value = -15366.286
$ws.Range("H96").Value = 3630.0833
$ws.Range("I96").Value = 3151.5715
$ws.Range("K96").Value = 3151.5715
$ws.Range("M96").Value = -1778.5715
